$d = $word.ActiveDocument

# Locate the first paragraph (ends with "...to help your imagine soar! ")
# and insert a brand-new empty paragraph right after it; the new paragraph
# will hold the Trial Mode text being added by this change.
$p1 = $d.Paragraphs(1)
$r1 = $p1.Range
$r1.Collapse(0)
$r1.InsertParagraphAfter() | Out-Null

# The freshly-created (empty) paragraph is now Paragraph 2; replace its
# full range (including the paragraph mark / placeholder run) with the real
# content so no stray empty run is left behind.
$p2 = $d.Paragraphs(2)
$target = $d.Range($p2.Range.Start, $p2.Range.End)

# Build the paragraph contents (runs + proofErr spell-check markers) as raw
# WordprocessingML so the run/proofErr boundaries match exactly.
$ellipsis = [char]0x2026
$paraXml = '<w:p>' +
  '<w:r><w:t xml:space="preserve">Just push the big button to get your next suggestion. Change lists by swiping to the left or right, or click the ellipse (' + $ellipsis + ') on the application bar to get a list of lists. </w:t></w:r>' +
  '<w:proofErr w:type="spellStart"/><w:r><w:t>Improv</w:t></w:r><w:proofErr w:type="spellEnd"/>' +
  '<w:r><w:t xml:space="preserve"> </w:t></w:r>' +
  '<w:proofErr w:type="spellStart"/><w:r><w:t>Suggester</w:t></w:r><w:proofErr w:type="spellEnd"/>' +
  '<w:r><w:t xml:space="preserve"> has a simple non-distracting look so you can concentrate in the suggestion; </w:t></w:r>' +
  '<w:r><w:t>i</w:t></w:r>' +
  '<w:r><w:t>t honors yo</w:t></w:r>' +
  '<w:r><w:t>u theme since you know what colors you like better than we do.</w:t></w:r>' +
  '</w:p>'

$packageXml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
  '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
  '<pkg:xmlData>' +
  '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
  '<w:body>' + $paraXml + '</w:body>' +
  '</w:document>' +
  '</pkg:xmlData></pkg:part></pkg:package>'

$target.InsertXML($packageXml) | Out-Null
